# "10 years Finalization data"
#
# The raw monthly sheet ("Data Harian - Table") mixes station metadata
# (rows 1-5), the actual daily-observations table (rows 9-40: one header
# row + 31 days of July data) and footnotes (rows 44-56).
#
# The finalization step extracts just the clean data table (header +
# daily rows) onto a brand-new "Sheet1" placed right after the raw
# sheet, and that new sheet becomes the active/selected tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Match the source sheet's final selection state (set before the new
# sheet is created/activated, so sheet1 doesn't end up "tabSelected").
$ws1.Range("A9:K40").Select()

# Insert the new sheet right after the raw data sheet.
$new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)

# Copy the header row + the 31 daily rows (values, number formats and
# styles) into the new sheet starting at A1.
$ws1.Range("A9:K40").Copy($new.Range("A1"))

# Select the whole copied range and make this the active sheet/tab.
$new.Range("A1:K32").Select()
$new.Activate()
